# Update market/profit data cells per scheduled runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 833.3333
$ws.Range("I2").Value = 500
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 500
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = -387
$ws.Range("N2").Value = -1226
$ws.Range("H40").Value = 1829.0834
$ws.Range("I40").Value = 1683.5
$ws.Range("J40").Value = 1877.6111
$ws.Range("K40").Value = 1683.5
$ws.Range("L40").Value = 1877.6111
$ws.Range("M40").Value = -1508.5
$ws.Range("N40").Value = -2227.6111
$ws.Range("H48").Value = 1200
$ws.Range("I48").Value = 300
$ws.Range("K48").Value = 900
$ws.Range("M48").Value = -608
$ws.Range("H53").Value = 1305.409
$ws.Range("I53").Value = 2061.182
$ws.Range("J53").Value = 549.63635
$ws.Range("K53").Value = 2061.182
$ws.Range("L53").Value = 549.63635
$ws.Range("M53").Value = -1424.182
$ws.Range("N53").Value = -1823.63635
$ws.Range("H56").Value = 1200
$ws.Range("I56").Value = 300
$ws.Range("K56").Value = 900
$ws.Range("M56").Value = -366
$ws.Range("H64").Value = 13891955
$ws.Range("I64").Value = 13891955
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 13891955
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -13891707
$ws.Range("H67").Value = 13891955
$ws.Range("I67").Value = 13891955
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 13891955
$ws.Range("L67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -13891097
$ws.Range("H98").Value = 2238.8462
$ws.Range("I98").Value = 2734.1667
$ws.Range("J98").Value = 1814.2858
$ws.Range("K98").Value = 2734.1667
$ws.Range("L98").Value = 1814.2858
$ws.Range("M98").Value = -1236.1667
$ws.Range("N98").Value = -4810.2858
$ws.Range("H100").Value = 2913.2666
$ws.Range("I100").Value = 2480
$ws.Range("J100").Value = 3779.8
$ws.Range("K100").Value = 2480
$ws.Range("L100").Value = 3779.8
$ws.Range("M100").Value = -1939
$ws.Range("N100").Value = -4861.8
$ws.Range("H122").Value = 2238.8462
$ws.Range("I122").Value = 2734.1667
$ws.Range("J122").Value = 1814.2858
$ws.Range("K122").Value = 8202.500100000001
$ws.Range("L122").Value = 5442.857400000001
$ws.Range("M122").Value = -5752.500100000001
$ws.Range("N122").Value = -10342.8574
$ws.Range("H129").Value = 1473.8889
$ws.Range("J129").Value = 1479.3258
$ws.Range("L129").Value = 4437.9774
$ws.Range("N129").Value = -14437.9774
$ws.Range("H138").Value = 8931570
$ws.Range("I138").Value = 1415.3334
$ws.Range("J138").Value = 62512500
$ws.Range("K138").Value = 4246.0002
$ws.Range("L138").Value = 187537500
$ws.Range("M138").Value = 893.9997999999996
$ws.Range("N138").Value = -187547780

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1232.5238
$ws.Range("I2").Value = 870
$ws.Range("K2").Value = 870
$ws.Range("M2").Value = -757
$ws.Range("H32").Value = 9946.714
$ws.Range("I32").Value = 11030.419
$ws.Range("J32").Value = 6362.154
$ws.Range("K32").Value = 11030.419
$ws.Range("L32").Value = 6362.154
$ws.Range("M32").Value = -10743.419
$ws.Range("N32").Value = -6936.154
$ws.Range("H45").Value = 3526.5
$ws.Range("I45").Value = 3942.4
$ws.Range("J45").Value = 2833.3333
$ws.Range("K45").Value = 3942.4
$ws.Range("L45").Value = 2833.3333
$ws.Range("M45").Value = -3565.4
$ws.Range("N45").Value = -3587.3333
$ws.Range("H116").Value = 1232.5238
$ws.Range("I116").Value = 870
$ws.Range("K116").Value = 870
$ws.Range("M116").Value = 1424

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1232.5238
$ws.Range("I3").Value = 870
$ws.Range("K3").Value = 870
$ws.Range("M3").Value = -756
$ws.Range("H20").Value = 827.2857
$ws.Range("I20").Value = 942.26666
$ws.Range("J20").Value = 539.8333
$ws.Range("K20").Value = 942.26666
$ws.Range("L20").Value = 539.8333
$ws.Range("M20").Value = -695.26666
$ws.Range("N20").Value = -1033.8333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7940607
$ws.Range("I31").Value = 7572.263
$ws.Range("J31").Value = 14493984
$ws.Range("K31").Value = 7572.263
$ws.Range("L31").Value = 14493984
$ws.Range("M31").Value = -7277.263
$ws.Range("N31").Value = -14494574
$ws.Range("H34").Value = 7940607
$ws.Range("I34").Value = 7572.263
$ws.Range("J34").Value = 14493984
$ws.Range("K34").Value = 7572.263
$ws.Range("L34").Value = 14493984
$ws.Range("M34").Value = -7370.263
$ws.Range("N34").Value = -14494388
$ws.Range("H107").Value = 399.11765
$ws.Range("I107").Value = 257.69565
$ws.Range("J107").Value = 694.8182
$ws.Range("K107").Value = 257.69565
$ws.Range("L107").Value = 694.8182
$ws.Range("M107").Value = 1662.30435
$ws.Range("N107").Value = -4534.8182

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 7411341.5
$ws.Range("I122").Value = 13337074
$ws.Range("J122").Value = 4175
$ws.Range("K122").Value = 40011222
$ws.Range("L122").Value = 12525
$ws.Range("M122").Value = -40008772
$ws.Range("N122").Value = -17425

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2363.0588
$ws.Range("I16").Value = 2363.0588
$ws.Range("K16").Value = 2363.0588
$ws.Range("M16").Value = -2193.0588
$ws.Range("H46").Value = 874
$ws.Range("I46").Value = 826.1539
$ws.Range("K46").Value = 826.1539
$ws.Range("M46").Value = -638.1539
$ws.Range("H55").Value = 483.58334
$ws.Range("I55").Value = 203.16667
$ws.Range("J55").Value = 764
$ws.Range("K55").Value = 203.16667
$ws.Range("L55").Value = 764
$ws.Range("M55").Value = -30.16667000000001
$ws.Range("N55").Value = -1110
$ws.Range("H68").Value = 1466.25
$ws.Range("I68").Value = 1150
$ws.Range("J68").Value = 1571.6666
$ws.Range("K68").Value = 1150
$ws.Range("L68").Value = 1571.6666
$ws.Range("M68").Value = -401
$ws.Range("N68").Value = -3069.6666
$ws.Range("H71").Value = 1466.25
$ws.Range("I71").Value = 1150
$ws.Range("J71").Value = 1571.6666
$ws.Range("K71").Value = 5750
$ws.Range("L71").Value = 7858.333000000001
$ws.Range("M71").Value = -2006
$ws.Range("N71").Value = -15346.333
$ws.Range("H122").Value = 6156.8335
$ws.Range("I122").Value = 10457.714
$ws.Range("J122").Value = 4385.8823
$ws.Range("K122").Value = 31373.142
$ws.Range("L122").Value = 13157.6469
$ws.Range("M122").Value = -28923.142
$ws.Range("N122").Value = -18057.6469
